$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell G1: "Material Type"
$hdr = $ws.Range("G1")
$hdr.Value = "Material Type"

# Start from the existing bold header style (F1, cellXf 4: fontId=1 fillId=3)
# so the font/fill already carry Bold + center alignment, then tweak the
# color of both to explicit white-on-black in a single write each - this
# keeps the style tables from growing more entries than necessary.
$srcHeader = $ws.Range("F1")
$srcHeader.Copy()
$hdr.PasteSpecial(-4122)
$hdr.Font.Color = 16777215
$hdr.Interior.Color = 0
$hdr.Interior.PatternColor = 0

# Data column G2:G24: "DNA:Genomic", using the same centered style as the
# existing C column values (cellXf 3).
$dataRange = $ws.Range("G2:G24")
$dataRange.Value = "DNA:Genomic"

$srcData = $ws.Range("C2")
$srcData.Copy()
$dataRange.PasteSpecial(-4122)
$dataRange.Value = "DNA:Genomic"

# Header selection now sits on the newly-added column.
$ws.Range("G1:G24").Select()
